# NaCl NIOSH Test Results - add a "Problems" column before "Comments" on
# the header row (row 5): E5 becomes "Problems" and the existing
# "Comments" header moves from G5 into F5 (G5 is cleared), shrinking the
# sheet's used range from column G to column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: reading back `.Value` on this host returns a reflection-info
# string rather than the cell's contents, so use `.Value2` for reads.
$commentsHeader = $ws.Range("G5").Value2

$ws.Range("E5").Value = "Problems"
$ws.Range("F5").Value = $commentsHeader
$ws.Range("G5").Clear()

# Match the bold header formatting used by the rest of row 5.
$ws.Range("E5").Font.Bold = $true
$ws.Range("F5").Font.Bold = $true

# Restore the view captured in the saved workbook: cell F6 selected and
# zoomed in to 110%.
$ws.Range("F6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 110
$null
